$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Title
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Quantum Entanglement: Unveiling the Mysteries", $false, $false, $false,
    $false, $false, $true, 1, $false,
    "Unlocking the Secrets of Life: An Exploration of Biology for High School Students",
    2)

# ------------------------------------------------------------------
# 2. Author name
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    " Alex Rodriguez", $false, $false, $false, $false, $false, $true, 1,
    $false, " Olivia Brown", 2)

# ------------------------------------------------------------------
# 3. Email paragraph -> "at"
# ------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$rng3 = $d.Range($p3.Range.Start, $p3.Range.End)
$rng3.Text = "at"

# ------------------------------------------------------------------
# 4. Body paragraph: rewrite the three quantum-physics blocks into
#    three biology blocks, each with a numbered heading followed by a
#    manual line break, while keeping the existing double line breaks
#    that separate the three blocks.
# ------------------------------------------------------------------
$oldBlock1 = "In the vast and enigmatic realm of quantum physics, the concept of quantum entanglement holds a prominent place, captivating the minds of scientists and philosophers alike. This phenomenon, unveiled in the 1930s, defies classical intuition and unveils a world of interconnectedness that transcends the confines of space and time. In this essay, we will delve into the complexities of quantum entanglement, exploring its profound implications on our understanding of reality and its potential applications in various fields."
$newBlock1 = "1. Journey into the Realm of Living Organisms:^lEnter the captivating world of biology, a realm that unravels the mysteries of life. Investigate the intricacies of organisms, their structures, and their functions. Unravel the enigmas of cellular biology, genetics, evolution, and ecology. Dissect the miraculous tapestry of life, revealing the interconnectedness of all living beings and their interdependence with the environment."
$d.Content.Find.Execute($oldBlock1, $false, $false, $false, $false, $false, $true, 1, $false, $newBlock1, 2)

$oldBlock2 = "Quantum entanglement, in its essence, describes a peculiar correlation between two particles, wherein the properties of one particle instantaneously influence the properties of the other, irrespective of the distance separating them. This phenomenon appears to transcend the limits of causality, challenging our conventional notions of locality and causality. Experiments conducted over vast distances have repeatedly confirmed the existence of entanglement, highlighting its intrinsic nonlocal character."
$newBlock2 = "2. Unraveling the Secrets of Life:^lUnveil the secrets of life, exploring the fundamentals of biochemistry, genetics, and molecular biology. Delve into the fascinating world of DNA, the blueprint of life, and uncover the mechanisms of protein synthesis. Witness the elegant symphony of cellular processes, unraveling the mysteries of cell division, energy production, and metabolism. Investigate the remarkable adaptations and diversity of organisms, revealing the breathtaking beauty of nature's artistry."
$d.Content.Find.Execute($oldBlock2, $false, $false, $false, $false, $false, $true, 1, $false, $newBlock2, 2)

$oldBlock3 = "The implications of quantum entanglement extend far beyond the theoretical realm, potentially revolutionizing fields such as cryptography, computing, and communication. Quantum cryptography exploits the inherent randomness of entangled particles to create unbreakable codes, ensuring the secure transmission of information. Quantum computers, leveraging the superposition and entanglement of quantum bits, promise exponential speed-ups in certain computations, potentially transforming industries and driving groundbreaking discoveries. Quantum communication networks, utilizing entangled particles as carriers of information, hold the promise of ultra-fast, secure, and long-distance communication."
$newBlock3 = "3. Exploring Evolution and Ecology:^lEmbark on a voyage through the annals of evolution, tracing the remarkable journey of life's transformation over billions of years. Delve into the mechanisms of natural selection, genetic variation, and adaptation, understanding how organisms evolve to survive and thrive in their ever-changing environments. Dive into the intricate web of ecology, exploring the dynamic interactions between organisms and their ecosystems. Discover the delicate balance of nature, highlighting the interdependence of species and the crucial role of biodiversity."
$d.Content.Find.Execute($oldBlock3, $false, $false, $false, $false, $false, $true, 1, $false, $newBlock3, 2)

# ------------------------------------------------------------------
# 5. Summary paragraph
# ------------------------------------------------------------------
$oldSummary = "Quantum entanglement, defying classical intuition, reveals a realm of interconnectedness and nonlocality in the quantum world. With profound implications for our understanding of reality, it holds potential applications in cryptography, computing, and communication. The entanglement of particles serves as a foundation for secure codes, exponential speed-ups in computations, and ultra-fast communication networks. As we continue to unravel the intricacies of entanglement, we may unlock new avenues for technological advancements and gain deeper insights into the fundamental fabric of the universe."
$newSummary = "Biology, an awe-inspiring journey into the realm of life, unveils the profound secrets of living organisms. It explores the intricate mechanisms of cells, unravels the mysteries of inheritance and evolution, and navigates the interconnectedness of organisms and ecosystems. Biology empowers us to understand the intricacies of our own existence and the incredible diversity of life around us, fostering a profound appreciation for the natural world."
$d.Content.Find.Execute($oldSummary, $false, $false, $false, $false, $false, $true, 1, $false, $newSummary, 2)

# ------------------------------------------------------------------
# 6. Add a new empty paragraph at the very end of the document.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

# ------------------------------------------------------------------
# 7. Font fix-up: TimesNewToman -> Times New Roman everywhere.
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $d.Paragraphs($i).Range.Font.Name = "Times New Roman"
}
